$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values would otherwise be auto-converted to numbers by Excel
# are explicitly formatted as Text first, to preserve them as text strings (matching the source data).
$ws.Range("D2").Value = "66.797.87"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.602.85"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.55"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.72"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "2.600.51"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.19"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.10"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "3.067.28"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "66.895.71"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "2.598.42"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.55"
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.77"
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.86"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.32"
$ws.Range("E25").Value = "  -6.87%  "
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.15"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "2.735.71"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "539.45"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.58"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.25"
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.80"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("D46").Value = "0.0₆0289"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.17"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.572"
$ws.Range("E48").Value = "  -3.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("E51").Value = "  -1.92%  "
